$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.966.87"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.928.42"
$ws.Range("E3").Value = "  -1.73%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "373.59"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "101.11"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.10%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("E8").Value = "  +0.00%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.579"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.84%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.23"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "3.396.00"
$ws.Range("E13").Value = "  -1.51%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "17.86"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.30"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "2.936.17"
$ws.Range("E16").Value = "  -1.38%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.971"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "50.920.62"
$ws.Range("E18").Value = "  -1.09%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.69%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.47"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.97%  "
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -1.22%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "263.55"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "68.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.42%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +6.69%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.99"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +6.45%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.85%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.61%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "25.53"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.86%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.80"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.63%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "50.94"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0450"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "33.36"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.50%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.33%  "
$ws.Range("E37").Value = "  -0.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.49%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.52"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "16.27"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.67%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "120.71"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "20.87"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.57%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -3.15%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.19"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "1.964.86"
$ws.Range("E49").Value = "  -3.58%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0341"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.02%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "5.01"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.11%  "
